# Updates the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# tracker sheet with freshly scraped figures.
#
# Some new Price figures (e.g. "0.998", "8.03") are valid numeric literals.
# A plain `Range.Value = "0.998"` assignment would make Excel store them as
# numbers, but the sheet represents prices as plain text (to preserve
# formatting such as the thousand-separator dots used for BTC/ETH, e.g.
# "58.600.26"). To force text entry for those cells we prefix the literal
# with an apostrophe (Excel's "treat as text" entry marker) and then strip
# the resulting quote-prefix cell format with ClearFormats() so the saved
# cell carries no extra style, matching the rest of the column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.600.26'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '3.151.99'
$ws.Range("E3").Value = '  +1.07%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''531.62'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").Value = '''139.76'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.20%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.520'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +12.31%  '
$ws.Range("D9").Value = '''7.34'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.67%  '
$ws.Range("D10").Value = '''0.426'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.17%  '
$ws.Range("E11").Value = '  +2.30%  '
$ws.Range("E12").Value = '  +2.71%  '
$ws.Range("D13").Value = '3.690.76'
$ws.Range("E13").Value = '  +0.87%  '
$ws.Range("D14").Value = '''25.74'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.08%  '
$ws.Range("E15").Value = '  +4.78%  '
$ws.Range("D16").Value = '58.623.30'
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("D17").Value = '3.152.87'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("E18").Value = '  +3.56%  '
$ws.Range("D19").Value = '''13.01'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.36%  '
$ws.Range("D20").Value = '''8.13'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.75%  '
$ws.Range("D21").Value = '''371.22'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +5.93%  '
$ws.Range("D22").Value = '''5.81'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.92%  '
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = '''69.85'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.34%  '
$ws.Range("D25").Value = '''0.516'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +2.56%  '
$ws.Range("D26").Value = '''0.167'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").Value = '''0.998'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '''8.03'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +11.37%  '
$ws.Range("D29").Value = '0.0₃0861'
$ws.Range("E29").Value = '  -0.97%  '
$ws.Range("D30").Value = '''1.88'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").Value = '''6.10'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.93%  '
$ws.Range("D32").Value = '''21.84'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.94%  '
$ws.Range("E33").Value = '  +5.36%  '
$ws.Range("E34").Value = '  +1.49%  '
$ws.Range("D35").Value = '''159.34'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").Value = '''6.27'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.94%  '
$ws.Range("E37").Value = '  +7.89%  '
$ws.Range("D38").Value = '''25.24'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.65%  '
$ws.Range("D39").Value = '2.658.19'
$ws.Range("E39").Value = '  +10.93%  '
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("D41").Value = '''0.0681'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.86%  '
$ws.Range("D42").Value = '''4.20'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.78%  '
$ws.Range("D43").Value = '''0.710'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.95%  '
$ws.Range("D44").Value = '''38.61'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +4.15%  '
$ws.Range("D45").Value = '''0.0283'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +6.62%  '
$ws.Range("D46").Value = '''0.999'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").Value = '3.191.62'
$ws.Range("E47").Value = '  +0.95%  '
$ws.Range("E48").Value = '  +12.06%  '
$ws.Range("D49").Value = '''0.983'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +2.23%  '
$ws.Range("E50").Value = '  +2.53%  '
$ws.Range("D51").Value = '''20.07'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.46%  '
